$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Drop the stray _GoBack bookmark that currently sits in the title
#    paragraph ("LIST OF ABBREVIATIONS"). It is hidden from
#    Bookmarks.Count (like real Word), but is still addressable by name.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $null = $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Add a new "CS" entry right after the "AS" (added subclass
#    (restriction)) line, consolidating the _GoBack bookmark onto the
#    end of this freshly added paragraph.
# ---------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^AS\t") {
        $target = $p
    }
}

$null = $target.Range.InsertParagraphAfter()
$newPara = $target.Next()

$fragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:r><w:t>CS</w:t></w:r>' +
            '<w:r><w:tab/><w:t>changed subclass (restriction)</w:t></w:r>' +
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
            '<w:bookmarkEnd w:id="0"/>' +
            '</w:p>'

$null = $newPara.Range.InsertXML($fragment)

Write-Output ("New paragraph text: [" + $newPara.Range.Text + "]")
Write-Output ("_GoBack present: " + $d.Bookmarks.Exists("_GoBack"))
